$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '94.085.95'
$ws.Range("E2").Value = '  -2.97%  '

Set-TextValue $ws.Range("D3") '3.435.38'
$ws.Range("E3").Value = '  +2.73%  '

$ws.Range("E4").Value = '  +0.11%  '

Set-TextValue $ws.Range("D5") '236.70'
$ws.Range("E5").Value = '  -5.60%  '

Set-TextValue $ws.Range("D6") '635.37'
$ws.Range("E6").Value = '  -3.39%  '

Set-TextValue $ws.Range("D7") '1.42'
$ws.Range("E7").Value = '  -0.25%  '

Set-TextValue $ws.Range("D8") '0.394'
$ws.Range("E8").Value = '  -7.32%  '

$ws.Range("E9").Value = '  +0.12%  '

Set-TextValue $ws.Range("D10") '0.966'
$ws.Range("E10").Value = '  -5.17%  '

Set-TextValue $ws.Range("D11") '3.437.44'
$ws.Range("E11").Value = '  +2.87%  '

Set-TextValue $ws.Range("D12") '41.86'
$ws.Range("E12").Value = '  +2.00%  '

Set-TextValue $ws.Range("D13") '0.196'
$ws.Range("E13").Value = '  -5.09%  '

Set-TextValue $ws.Range("D14") '6.14'
$ws.Range("E14").Value = '  +0.48%  '

Set-TextValue $ws.Range("D15") '94.074.19'
$ws.Range("E15").Value = '  -2.74%  '

Set-TextValue $ws.Range("D16") '4.078.11'
$ws.Range("E16").Value = '  +2.79%  '

Set-TextValue $ws.Range("D17") '0.0000251'
$ws.Range("E17").Value = '  -0.91%  '

Set-TextValue $ws.Range("D18") '8.34'
$ws.Range("E18").Value = '  -5.54%  '

Set-TextValue $ws.Range("D19") '3.433.34'
$ws.Range("E19").Value = '  +3.40%  '

Set-TextValue $ws.Range("D20") '17.53'
$ws.Range("E20").Value = '  -0.70%  '

Set-TextValue $ws.Range("D21") '11.26'
$ws.Range("E21").Value = '  +4.91%  '

Set-TextValue $ws.Range("D22") '0.497'
$ws.Range("E22").Value = '  -10.98%  '

Set-TextValue $ws.Range("D23") '495.00'
$ws.Range("E23").Value = '  -2.98%  '

Set-TextValue $ws.Range("D24") '3.12'
$ws.Range("E24").Value = '  -6.42%  '

Set-TextValue $ws.Range("D25") '6.55'
$ws.Range("E25").Value = '  -0.86%  '

Set-TextValue $ws.Range("D26") '0.0000188'
$ws.Range("E26").Value = '  -5.71%  '

Set-TextValue $ws.Range("D27") '90.76'
$ws.Range("E27").Value = '  -6.30%  '

Set-TextValue $ws.Range("D28") '3.618.94'
$ws.Range("E28").Value = '  +2.84%  '

Set-TextValue $ws.Range("D29") '11.88'
$ws.Range("E29").Value = '  -2.41%  '

Set-TextValue $ws.Range("D30") '11.65'
$ws.Range("E30").Value = '  +1.36%  '

Set-TextValue $ws.Range("D31") '1.00'
$ws.Range("E31").Value = '  -0.02%  '

Set-TextValue $ws.Range("D32") '2.72'
$ws.Range("E32").Value = '  +6.72%  '

Set-TextValue $ws.Range("D33") '0.134'
$ws.Range("E33").Value = '  -8.01%  '

Set-TextValue $ws.Range("D34") '0.180'
$ws.Range("E34").Value = '  -4.56%  '

Set-TextValue $ws.Range("D35") '0.999'
$ws.Range("E35").Value = '  -0.02%  '

Set-TextValue $ws.Range("D36") '30.15'
$ws.Range("E36").Value = '  +5.88%  '

Set-TextValue $ws.Range("D37") '0.560'
$ws.Range("E37").Value = '  +0.42%  '

Set-TextValue $ws.Range("D38") '535.61'
$ws.Range("E38").Value = '  +5.40%  '

Set-TextValue $ws.Range("D39") '7.59'
$ws.Range("E39").Value = '  -3.26%  '

Set-TextValue $ws.Range("D40") '1.43'
$ws.Range("E40").Value = '  -4.49%  '

$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D41") '1.00'
$ws.Range("E41").Value = '  -0.12%  '

Set-TextValue $ws.Range("D42") '0.923'
$ws.Range("E42").Value = '  +9.84%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D43") '0.150'
$ws.Range("E43").Value = '  -1.14%  '

Set-TextValue $ws.Range("D44") '24.03'
$ws.Range("E44").Value = '  -1.43%  '

Set-TextValue $ws.Range("D45") '1.67'
$ws.Range("E45").Value = '  -0.81%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D46") '0.0408'
$ws.Range("E46").Value = '  -6.72%  '

$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D47") '5.51'
$ws.Range("E47").Value = '  -2.88%  '

Set-TextValue $ws.Range("D48") '3.48'
$ws.Range("E48").Value = '  -4.70%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D49") '53.52'
$ws.Range("E49").Value = '  -2.18%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D50") '2.13'
$ws.Range("E50").Value = '  +6.89%  '

Set-TextValue $ws.Range("D51") '3.17'
$ws.Range("E51").Value = '  +1.67%  '
